# Insert a new weekly record at row 87, shifting all existing records
# (previously rows 87-208) down by one row to rows 88-209.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new weekly entry.
$ws.Range("A87").Value2 = 8
$ws.Range("B87").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C87").Value2 = "Coquimbo"
$ws.Range("D87").Value2 = 44571
$ws.Range("E87").Value2 = 4
$ws.Range("F87").Value2 = 100112012
$ws.Range("G87").Value2 = "Espinaca"
$ws.Range("H87").Value2 = "Sin especificar"
$ws.Range("I87").Value2 = "Primera"
$ws.Range("J87").Value2 = 2400
$ws.Range("K87").Value2 = 400
$ws.Range("L87").Value2 = 500
$ws.Range("M87").Value2 = 450
$ws.Range("N87").Value2 = "`$/atado 300 a 500 gramos"
$ws.Range("O87").Value2 = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P87").Value2 = 900
$ws.Range("Q87").Value2 = 0.5
$ws.Range("R87").Value2 = "Hortaliza"

# Note: the date number format for column D (numFmtId 165) on row 87 is
# already inherited correctly from the Insert() operation above, so no
# further style assignment is necessary here.
